# Generate Report for Handback
#
# This updates the localization-status workbook to reflect a failed
# handback transform for the 0348d7f1-... file:
#   - Status changes from "Ready for handoff" to "Handback transform failed"
#     on every sheet that shows it (Overview, zh-cn, de-de).
#   - The "Error Detail" column (P) on the zh-cn and de-de sheets gets a
#     new diagnostic message for that row, and the column is widened so
#     the message is readable.

$wb = $excel.ActiveWorkbook

$ov = $wb.Worksheets.Item("Overview")
$zh = $wb.Worksheets.Item("zh-cn")
$de = $wb.Worksheets.Item("de-de")

# Update the Status text everywhere it appears for the 0348d7f1 row (row 3
# on every sheet) so they all keep sharing the same string.
[void]$ov.Cells.Replace("Ready for handoff", "Handback transform failed")
[void]$zh.Cells.Replace("Ready for handoff", "Handback transform failed")
[void]$de.Cells.Replace("Ready for handoff", "Handback transform failed")

# New Error Detail messages for that row.
$zh.Range("P3").Value = "Handback file name: oeny3vvr.ph1 is different with handoff file name: 0348d7f1-06c9-47b5-83b4-f39cbabd09a0.2b5b18774d0b441e0db3b0a22239aa1cb20def2d.zh-cn."
$de.Range("P3").Value = "Handback file name: oeny3vvr.ph1 is different with handoff file name: 0348d7f1-06c9-47b5-83b4-f39cbabd09a0.2b5b18774d0b441e0db3b0a22239aa1cb20def2d.de-de."

# Widen the Error Detail column (P) on both sheets to fit the new text.
# Excel's ColumnWidth object-model units run 5/6 of a character above the
# raw OOXML column width, so subtract that offset to land on width=40.
$zh.Range("P1").ColumnWidth = 39.166666666666664
$de.Range("P1").ColumnWidth = 39.166666666666664
